$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201, shifting existing rows 201-215 down to 202-216
$ws.Rows.Item(201).Insert()

# Populate the new row 201 with data (copy template values from row 202, then override changed cells)
$ws.Range("A201").Value = 5
$ws.Range("B201").Value = "Macroferia Regional de Talca"
$ws.Range("C201").Value = "Maule"
$ws.Range("D201").Value = 44783
$ws.Range("D201").NumberFormat = $ws.Range("D202").NumberFormat
$ws.Range("E201").Value = 7
$ws.Range("F201").Value = 100112017
$ws.Range("G201").Value = "Apio"
$ws.Range("H201").Value = "Americana (o)"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 700
$ws.Range("K201").Value = 10000
$ws.Range("L201").Value = 10000
$ws.Range("M201").Value = 10000
$ws.Range("N201").Value = "$/docena de matas"
$ws.Range("O201").Value = "Provincia del Elquí"
$ws.Range("P201").Value = 1667
$ws.Range("Q201").Value = 6
$ws.Range("R201").Value = "Hortaliza"
